# Add carjacking data for 2021-12-04 through 2021-12-12.
# Updates the sheet title, the running "through" label, and the day-level
# counts (column B = December 2021 totals, N/Z/AX/BJ/BV etc = day-of-month
# detail columns) for the neighborhoods that saw new incidents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet name reflects the new "through" date.
$ws.Name = "Through 2021-12-04"

# Header cell with the same "through" date, spelled out.
$ws.Range("B1").Value = "December 2021 (through December 04)"

# --- Row 3: North Lawndale ---
$ws.Range("N3").Value = 1

# --- Row 6: Garfield Park ---
$ws.Range("N6").Value = 3

# --- Row 7: Austin ---
$ws.Range("B7").Value = 3
$ws.Range("BJ7").Value = 1

# --- Row 8: Grand Crossing ---
$ws.Range("B8").Value = 3
$ws.Range("Z8").Value = 1

# --- Row 9: Humboldt Park ---
$ws.Range("Z9").Value = 1

# --- Row 10: Douglas ---
$ws.Range("B10").Value = 3
$ws.Range("AX10").Value = 2

# --- Row 11: Roseland ---
$ws.Range("AX11").Value = 1

# --- Row 17: Bridgeport ---
$ws.Range("N17").Value = 1
$ws.Range("AX17").Value = 1

# --- Row 25: Wicker Park ---
$ws.Range("N25").Value = 1

# --- Row 29: Woodlawn ---
$ws.Range("B29").Value = 2
$ws.Range("AL29").Value = 1

# --- Row 33: Auburn Gresham ---
$ws.Range("AX33").Value = 2
$ws.Range("BV33").Value = 1

# --- Row 34: New City ---
$ws.Range("B34").Value = 1

# --- Row 35: Near South Side ---
$ws.Range("B35").Value = 1

# --- Row 37: Gage Park ---
$ws.Range("BJ37").Value = 1

# --- Row 38: Calumet Heights ---
$ws.Range("B38").Value = 1

# --- Row 40: Chinatown ---
$ws.Range("B40").Value = 1

# --- Row 43: Ukrainian Village ---
$ws.Range("B43").Value = 1

# --- Row 63: Portage Park ---
$ws.Range("AL63").Value = 2

# --- Row 67: Avondale ---
$ws.Range("B67").Value = 2

# --- Row 79: Irving Park ---
$ws.Range("N79").Value = 1
$ws.Range("AX79").Value = 2

# --- Row 94: Sauganash, Forest Glen ---
$ws.Range("N94").Value = 1
